# Weekly update: the dataset rows (D:R, rows 2-30) are reshuffled.
# Each destination row receives the full D:R content that used to live
# at a different source row. Column A, B and C are unchanged (constant
# across the whole sheet), so we only need to touch D:R.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> old row number whose D:R data it now holds.
$rowMap = @{
    2  = 13
    3  = 14
    4  = 20
    5  = 21
    6  = 22
    7  = 28
    8  = 29
    9  = 16
    10 = 19
    11 = 18
    12 = 6
    13 = 2
    14 = 23
    15 = 24
    16 = 10
    17 = 25
    18 = 30
    19 = 9
    20 = 4
    21 = 15
    22 = 26
    23 = 27
    24 = 5
    25 = 3
    26 = 17
    27 = 8
    28 = 7
    29 = 11
    30 = 12
}

# Snapshot every source cell's value (D:R, columns 4-18) before any writes,
# so overlapping reads/writes in the permutation don't clobber each other.
$colStart = 4   # D
$colEnd   = 18  # R

$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $ws.Cells.Item($newRow, $c).Value2 = $snapshot["$oldRow-$c"]
    }
}
